$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 values
$ws.Range("B2").Value = 4.737029407806765
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 1.2270924346332597
$ws.Range("E2").Value = 0.69081181766519206

# Row 3 values
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 6.0114943706668864
$ws.Range("D3").Value = 6.6860536061983291
$ws.Range("E3").Value = -1.4488257248292093

# Update selection to match new range
$ws.Range("B1:E3").Select()
